# Update the table style used by the "Data Sources" tables from the
# default LFX-exported style GUID to the new style GUID.
#
# Old style: {A1FA7028-C1A7-44A9-B84D-1D3E6534E97C}
# New style: {AA5A75A7-0268-4E23-BE90-BAA816196550}

$oldStyle = "{A1FA7028-C1A7-44A9-B84D-1D3E6534E97C}"
$newStyle = "{AA5A75A7-0268-4E23-BE90-BAA816196550}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyle) {
                $table.ApplyStyle($newStyle)
            }
        }
    }
}
